# Adds two new company rows (Project Prometheus / XtalPi) to the
# "AI_MSE_companies" tracker sheet, right after the existing Entalpic row.
#
# Columns: A=company  B=website  C=publicly announced funding
#          D=notable scientists  E=mission or domain  F=product / platform
#          G=sources (URLs)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 (Entalpic) uses style s=3 (vertical-top, wrap-text). Pre-copy that
# formatting onto every new cell that needs it. Copying from an
# already-populated donor cell doesn't introduce new shared strings (its
# text is already in the table), and the values below overwrite the
# copied text right after.
$ws.Range("A35").Copy($ws.Range("A36"))
$ws.Range("A35").Copy($ws.Range("C36"))
$ws.Range("A35").Copy($ws.Range("D36"))
$ws.Range("A35").Copy($ws.Range("A37"))
$ws.Range("A35").Copy($ws.Range("E37"))

# Fill in the cell values in the same order the author originally typed
# them, so new shared-string entries land in the same sequence.

# --- Row 36: Project Prometheus ---
$ws.Range("A36").Value = "Project Prometheus"
$ws.Range("C36").Value = "`$6.2B from Bezos fund"
$ws.Range("G36").Value = "https://www.nytimes.com/2025/11/17/technology/bezos-project-prometheus.html"
$ws.Range("D36").Value = "Jeff Bezos backing, Vik Bajaj (formerly Google X, Verily, Foresite Labs)"

# --- Row 37: XtalPi ---
$ws.Range("B37").Value = "https://en.xtalpi.com/"
$ws.Range("A37").Value = "XtalPi"
$ws.Range("E37").Value = "Drug discovery, lithium-sulfur batteries, "
$ws.Range("G37").Value = "https://www.texau.com/profiles/xtal-pi"
$ws.Range("C37").Value = "`$784M (`$400 million in Series D (August 11, 2021), `$268 million in Post-IPO (February 27, 2025), and `$ 116 million in Post-IPO (September 17, 2025).)"

# Wrapped-text row height for the newly filled row 36.
$ws.Rows.Item(36).RowHeight = 29

# Leave the selection where the author ended up after typing the new rows.
$ws.Range("C38").Select() | Out-Null
